# Daily attendance processing - 2026-01-20 16:49:07
# The "Recorded By" column (G) lists the users/systems that recorded a
# session. Previously the "System" token was written first, e.g.
#   "System, dnasr281@gmail.com"
# It now needs to list the human user first, then "System", e.g.
#   "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
